# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# D column holds price text like "68.524.69" or "0.0000190" that must stay
# literal strings (not be re-interpreted as numbers), so each write forces
# Text format first and then restores the Normal style afterwards (the
# source cells carry no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.524.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.723.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.75%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.60%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.721.79"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.146"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.366"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.35%  "

$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.215.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000190"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.522.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.721.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "373.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.63%  "

$ws.Range("E24").Value = "  +1.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.46%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  +3.77%  "

$ws.Range("E28").Value = "  +2.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000105"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "592.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.57%  "

$ws.Range("E31").Value = "  +0.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.04%  "

$ws.Range("E34").Value = "  +5.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.131"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.10%  "

$ws.Range("E40").Value = "  +2.53%  "

$ws.Range("E41").Value = "  +1.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.06%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.89%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0311"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.600"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "155.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.12%  "
